$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: replace computed date serials with plain text date labels
$ws.Range("A2").Value = "29/3/2020"

$ws.Range("A3").Formula = ""
$ws.Range("A3").Value = "30/3/2020"

$ws.Range("A4").Formula = ""
$ws.Range("A4").Value = "31/3/2020"

# Rows 5-9: drop the running-date formula, replace with literal (non-sequential) dates
$ws.Range("A5").Formula = ""
$ws.Range("A5").Value = 43834

$ws.Range("A6").Formula = ""
$ws.Range("A6").Value = 43865

$ws.Range("A7").Formula = ""
$ws.Range("A7").Value = 43894

$ws.Range("A8").Formula = ""
$ws.Range("A8").Value = 43925

$ws.Range("A9").Formula = ""
$ws.Range("A9").Value = 43955

# Row 10 becomes the new start of the running "+previous+1" formula chain;
# rows 11-51 (already holding "+previous+1" formulas) recalculate automatically.
$ws.Range("A10").Formula = "=+A9+1"

# Update the visible selection / scroll position recorded in the sheet view.
$ws.Range("A10").Select() | Out-Null
